# "with new goal of KTH map" - update the goal pose values on Sheet1
# and move the active selection to the newly-edited cell (B7), matching
# the intent of the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New goal values (msg.Pose.Position.Y / msg.Pose.Position.Z / msg.Pose.Orientation.X / msg.Pose.Orientation.Y)
$ws.Range("B2").Value = 255.29
$ws.Range("B3").Value = -514.31
$ws.Range("B7").Value = 0.481
$ws.Range("B8").Value = -0.877

# Move / leave the selection on B7 (matches the updated sheetView selection)
$ws.Range("B7").Select()
